$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows 2 (PLAYER_BULLET_DAMAGE), 3 (duplicate PLAYER_BULLET_SPEED),
# and 11 (PLAYER_HP). Delete from bottom to top so row numbers stay valid.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

$ws.Range("A7").Select()
